# Insert a new "abbr" column (state abbreviation) between the existing
# "state" column (C) and "most_populous_city" column (old D, now E).
# This shifts the old D:E columns (most_populous_city, city_pop) to E:F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (most_populous_city), shifting it
# and city_pop one column to the right.
$ws.Columns.Item(4).Insert()

# Header for the new column.
$ws.Range("D1").Value = "abbr"

# Two-letter USPS abbreviations, in the same row order (2-52) as the
# existing "state" column (alphabetical, with District of Columbia
# between Delaware and Florida).
$abbrs = @(
    "AL", "AK", "AZ", "AR", "CA", "CO", "CT", "DE", "DC", "FL",
    "GA", "HI", "ID", "IL", "IN", "IA", "KS", "KY", "LA", "ME",
    "MD", "MA", "MI", "MN", "MS", "MO", "MT", "NE", "NV", "NH",
    "NJ", "NM", "NY", "NC", "ND", "OH", "OK", "OR", "PA", "RI",
    "SC", "SD", "TN", "TX", "UT", "VT", "VA", "WA", "WV", "WI",
    "WY"
)

for ($i = 0; $i -lt $abbrs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $abbrs[$i]
}

# Re-apply the autofilter so its range grows from A1:E1 to A1:F1, and
# keep the workbook's hidden _FilterDatabase defined name in sync.
$ws.AutoFilterMode = $false | Out-Null
$ws.Range("A1:F1").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name() -eq "citypop!_FilterDatabase") {
        $n.RefersTo = "=citypop!`$A`$1:`$F`$1"
    }
}
